$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Update the password hash for the existing user (row 2 / M'mah Kamanda)
$ws.Range("D2").Value = '$2a$10$Bz4/5bkmPGFx.KNDLW2Us.iO2Q9dDQjt0wGkqqVFHyeLfHTjb.EF.'

# Add a new row for the new user M'mah Zombo
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'M''mah Zombo'
$ws.Range("C3").Value = 'zombo@gmail.com'
$ws.Range("D3").Value = '$2a$10$/sNF3JT9o2N3GlJj//AFE.bxtu9fT9CyTXrZD1iaVRM9g9nH8UaCa'
$ws.Range("E3").Value = 'agent'
$ws.Range("F3").Value = ''

$wb.Save()
